$wb = $excel.ActiveWorkbook

# --- Overview sheet: status summary columns mirror the per-locale status ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")

# Status column (B) now reflects a completed handback
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("B3").Value = "Handed back: in sync with en-US"

# Row 2 (b4425298...): fill in Latest Target File (E) and Latest Handback File (F)
$ws.Range("E2").Value = "b4425298-3a14-4a64-8351-b5a5aaffc974.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/bd7273533a9a44808b1a015d40160cc62fa78461/e2e/b4425298-3a14-4a64-8351-b5a5aaffc974.md", "", "", "b4425298-3a14-4a64-8351-b5a5aaffc974.md")

$ws.Range("F2").Value = "b4425298-3a14-4a64-8351-b5a5aaffc974.305f89c047dd2ee9fcf57d33bd0ac6eec4be2407.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b964e8a7cc7566f4abb86768de5b5a05abf59ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high/b4425298-3a14-4a64-8351-b5a5aaffc974.305f89c047dd2ee9fcf57d33bd0ac6eec4be2407.zh-cn.xlf", "", "", "b4425298-3a14-4a64-8351-b5a5aaffc974.305f89c047dd2ee9fcf57d33bd0ac6eec4be2407.zh-cn.xlf")

# Row 2: Latest Handback DateTime (G)
$ws.Range("G2").Value = "2016-03-07 07:06:11"

# Row 3 (d766af45...): fill in Latest Target File (E) and Latest Handback File (F)
$ws.Range("E3").Value = "d766af45-5c0d-4f2b-a6b5-c550a423a439.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/bd7273533a9a44808b1a015d40160cc62fa78461/e2e/d766af45-5c0d-4f2b-a6b5-c550a423a439.md", "", "", "d766af45-5c0d-4f2b-a6b5-c550a423a439.md")

$ws.Range("F3").Value = "d766af45-5c0d-4f2b-a6b5-c550a423a439.548865028afd694c0360e1c00cc07dde0e0d9da2.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b964e8a7cc7566f4abb86768de5b5a05abf59ba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high/d766af45-5c0d-4f2b-a6b5-c550a423a439.548865028afd694c0360e1c00cc07dde0e0d9da2.zh-cn.xlf", "", "", "d766af45-5c0d-4f2b-a6b5-c550a423a439.548865028afd694c0360e1c00cc07dde0e0d9da2.zh-cn.xlf")

# Row 3: Latest Handback DateTime (G)
$ws.Range("G3").Value = "2016-03-07 07:06:11"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("B3").Value = "Handed back: in sync with en-US"

$ws.Range("E2").Value = "b4425298-3a14-4a64-8351-b5a5aaffc974.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/bd7273533a9a44808b1a015d40160cc62fa78461/e2e/b4425298-3a14-4a64-8351-b5a5aaffc974.md", "", "", "b4425298-3a14-4a64-8351-b5a5aaffc974.md")

$ws.Range("F2").Value = "b4425298-3a14-4a64-8351-b5a5aaffc974.305f89c047dd2ee9fcf57d33bd0ac6eec4be2407.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/912b73bcbb49fdd3ae97b1936bbcbd6fef9e1ba1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high/b4425298-3a14-4a64-8351-b5a5aaffc974.305f89c047dd2ee9fcf57d33bd0ac6eec4be2407.de-de.xlf", "", "", "b4425298-3a14-4a64-8351-b5a5aaffc974.305f89c047dd2ee9fcf57d33bd0ac6eec4be2407.de-de.xlf")

$ws.Range("G2").Value = "2016-03-07 07:06:32"

$ws.Range("E3").Value = "d766af45-5c0d-4f2b-a6b5-c550a423a439.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/bd7273533a9a44808b1a015d40160cc62fa78461/e2e/d766af45-5c0d-4f2b-a6b5-c550a423a439.md", "", "", "d766af45-5c0d-4f2b-a6b5-c550a423a439.md")

$ws.Range("F3").Value = "d766af45-5c0d-4f2b-a6b5-c550a423a439.548865028afd694c0360e1c00cc07dde0e0d9da2.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/912b73bcbb49fdd3ae97b1936bbcbd6fef9e1ba1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high/d766af45-5c0d-4f2b-a6b5-c550a423a439.548865028afd694c0360e1c00cc07dde0e0d9da2.de-de.xlf", "", "", "d766af45-5c0d-4f2b-a6b5-c550a423a439.548865028afd694c0360e1c00cc07dde0e0d9da2.de-de.xlf")

$ws.Range("G3").Value = "2016-03-07 07:06:32"
